$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("RUNMANAGER")
$ws3 = $wb.Worksheets.Item("DATA")

# ---------------------------------------------------------------
# RUNMANAGER sheet: C2 execute flag flips from "yes" to "no"
# ---------------------------------------------------------------
$ws2.Cells.Item(2,3).Value = "no"

# New row 4 on RUNMANAGER (amazonDemoTest). Copy row3's formatting down first
# so the new cells inherit the same text-format styles (s=2 / s=3).
$ws2.Range("A3:E3").Copy()
$ws2.Range("A4:E4").PasteSpecial(-4122)

# Write B4 first so "Amazon Demo Test" becomes shared-string index 36
# (matches the order new strings were authored in the real edit).
$ws2.Cells.Item(4,2).Value = "Amazon Demo Test"

# ---------------------------------------------------------------
# DATA sheet: two new columns (G menuOption, H subMenuOption)
# ---------------------------------------------------------------

# Headers G1/H1 - copy format from F1 (same header style) first.
$ws3.Range("F1").Copy()
$ws3.Range("G1:H1").PasteSpecial(-4122)
$ws3.Cells.Item(1,7).Value = "menuOption"

# Data columns G2:H5 - copy format from an existing "text" styled cell
# (RUNMANAGER D2, which already carries the quotePrefix style s=3) then
# blank them out with an explicit quote-prefixed empty string so they end
# up as shared empty string literals just like the real workbook.
$ws2.Range("D2").Copy()
$ws3.Range("G2:H5").PasteSpecial(-4122)
$ws3.Cells.Item(2,7).Value = "'"

# New row 6 on DATA (amazonDemoTest row), replacing the old blank
# placeholder rows 6-15. Remove rows 7-15 first, shifting row6 up to stay,
# then populate row 6.
$ws3.Rows("7:15").Delete()

$ws3.Cells.Item(6,2).Value = "yes"
$ws3.Cells.Item(6,3).Value = "chrome"

# D6:F6 stay blank but quote-prefixed/text-styled like the rest of column.
$ws2.Range("D2:E2").Copy()
$ws3.Range("D6:E6").PasteSpecial(-4122)
$ws2.Range("E2").Copy()
$ws3.Range("F6").PasteSpecial(-4122)
$ws3.Range("D6:F6").Value = "'"

# H6 = Laptops (new string, before amazonDemoTest so indices line up)
$ws3.Cells.Item(6,8).Value = "Laptops"

# A4 on RUNMANAGER = amazonDemoTest (creates shared string, reused by DATA A6)
$ws2.Cells.Item(4,1).Value = "amazonDemoTest"
$ws3.Cells.Item(6,1).Value = "amazonDemoTest"

# Finish header H1 = subMenuOption
$ws3.Cells.Item(1,8).Value = "subMenuOption"

# G6 = "Mobiles, Computers"
$ws3.Cells.Item(6,7).Value = "Mobiles, Computers"

# ---------------------------------------------------------------
# RUNMANAGER row4 remaining cells: C4 yes, D4/E4 "1"
# ---------------------------------------------------------------
$ws2.Cells.Item(2,3).Copy()
$ws2.Cells.Item(4,3).PasteSpecial(-4163)
$ws2.Cells.Item(2,3).Value = "no"
$ws2.Cells.Item(4,3).Value = "yes"

$ws2.Cells.Item(2,4).Copy()
$ws2.Cells.Item(4,4).PasteSpecial(-4163)
$ws2.Cells.Item(2,5).Copy()
$ws2.Cells.Item(4,5).PasteSpecial(-4163)

# ---------------------------------------------------------------
# Column widths / selections / dimension are refreshed automatically,
# but nudge the "bestFit" columns wider to account for the longer text.
# ---------------------------------------------------------------
$ws2.Columns.Item(1).AutoFit()
$ws3.Columns.Item(1).AutoFit()
$ws3.Columns.Item(7).AutoFit()
$ws3.Columns.Item(8).AutoFit()

$ws2.Range("B10").Select()
$ws3.Range("G6").Select()
